$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L and M carry a "Text" (@) number format on these rows, so writing
# a numeric literal straight into .Value would get stored as a text string.
# Flip the cell to General just long enough to write the number, then put
# the original (text) number format back so the stored style is unchanged.
function Set-NumericValue($range, $value) {
    $fmt = $range.NumberFormat
    $range.NumberFormat = "general"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Update existing rows 357, 359, 360 (input values only; the B/H/J/K
#     formula cells recalc automatically) ---

# Row 357: Nb nouveaux cas positifs (C) 69 -> 68
$ws.Range("C357").Value = 68

# Row 359: Nb nouveaux cas positifs (C) 32 -> 46
#          Nb nouveaux deces a l'hopital (L) 2 -> 3
$ws.Range("C359").Value = 46
Set-NumericValue $ws.Range("L359") 3

# Row 360: Nb nouveaux cas positifs (C) 5 -> 50
#          Nb nouveaux deces a l'hopital (L) 0 -> 3
$ws.Range("C360").Value = 50
Set-NumericValue $ws.Range("L360") 3

# --- Fill in previously-blank rows 361-363 ---

# Row 361 (44247)
$ws.Range("C361").Value = 49
$ws.Range("E361").Value = 7
$ws.Range("F361").Value = 5
$ws.Range("G361").Value = 30
Set-NumericValue $ws.Range("L361") 0
Set-NumericValue $ws.Range("M361") 0

# Row 362 (44248)
$ws.Range("C362").Value = 27
$ws.Range("E362").Value = 9
$ws.Range("F362").Value = 6
$ws.Range("G362").Value = 31
Set-NumericValue $ws.Range("L362") 0
Set-NumericValue $ws.Range("M362") 0

# Row 363 (44249)
$ws.Range("C363").Value = 14
$ws.Range("E363").Value = 9
$ws.Range("F363").Value = 7
$ws.Range("G363").Value = 31
Set-NumericValue $ws.Range("L363") 0
Set-NumericValue $ws.Range("M363") 0

# --- Update the selected/active cell shown in the frozen-pane view ---
$ws.Range("S15").Select()
